$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.097.28"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "2.655.83"
$ws.Range("E3").Value = "  +1.95%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.599"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.63"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.76%  "
$ws.Range("E10").Value = "  +1.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.383"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.156"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("D13").Value = "3.126.29"
$ws.Range("E13").Value = "  +1.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +11.91%  "
$ws.Range("D15").Value = "61.076.41"
$ws.Range("E15").Value = "  +0.83%  "
$ws.Range("E16").Value = "  +1.23%  "
$ws.Range("D17").Value = "2.665.91"
$ws.Range("E17").Value = "  +1.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "352.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.529"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.39%  "
$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.162"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.94%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.79%  "
$ws.Range("D29").Value = "0.0₃0817"
$ws.Range("E29").Value = "  +3.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.93"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.55%  "
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "166.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.07"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +10.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.51"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.33"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.68"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.84%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "338.63"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +12.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.04"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.22%  "
$ws.Range("B40").Value = "SuiNetwork"
$ws.Range("C40").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.892"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.50%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "38.56"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.25"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.41"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "134.11"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.41%  "
$ws.Range("E45").Value = "  +1.44%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0249"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.10%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0562"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.67%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.616"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.68%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.998"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.26%  "
$ws.Range("D51").Value = "2.101.09"
$ws.Range("E51").Value = "  +3.63%  "
